$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The column-E "route" notes on rows 2-6 are being relocated to the bottom
# of column F (new rows 47-51), so clear the old E-column cells first.
$ws.Range("E2").Clear()
$ws.Range("E3").Clear()
$ws.Range("E4").Clear()
$ws.Range("E5").Clear()
$ws.Range("E6").Clear()

# Row 3 re-wraps to a shorter height once its column-E note is gone.
$ws.Range("B3").RowHeight = 45

# New admin-panel / site-map related todo items added in column B.
$ws.Range("B16").Value = "فوتر"
$ws.Range("B17").Value = "تصویر داریان در سوشال مدیا"
$ws.Range("B18").Value = "تگ ها درست نت"
$ws.Range("B19").Value = "تصاویر پنل ادمین"

# Re-append the relocated notes to the bottom of column F (new rows 47-51),
# in the same order they used to appear going down column E.
$ws.Range("F47").Value = "منیفست"

$ws.Range("F48").Value = "آخرین پست در بنر صفحه اصلی"
$ws.Range("F48").RowHeight = 60

$ws.Range("F49").Value = "لوگو اینستاگرام"
$ws.Range("F49").RowHeight = 30

$ws.Range("F50").Value = "save and continue"
$ws.Range("F50").RowHeight = 30

$ws.Range("F51").Value = "کتابخانه"

# Match the saved selection/scroll state.
$ws.Range("B20").Select()
